$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 344: new positive cases count 144 -> 146 ---
$ws.Range("C344").Value = 146

# --- Row 345: new positive cases count 93 -> 94 ---
$ws.Range("C345").Value = 94

# --- Row 346: new positive cases count 100 -> 101 ---
$ws.Range("C346").Value = 101

# --- Row 347: new extra-hospital death (M) 0 -> 1 (C347 unchanged) ---
# L/M columns are formatted as Text ("@") but hold numeric values, so
# temporarily switch to a numeric format while writing, then restore it.
$lm347 = $ws.Range("M347")
$lm347.NumberFormat = "General"
$lm347.Value = 1
$lm347.NumberFormat = "@"

# --- Row 348: new positive cases 47 -> 46; new extra-hospital death 0 -> 1 ---
$ws.Range("C348").Value = 46
$lm348 = $ws.Range("M348")
$lm348.NumberFormat = "General"
$lm348.Value = 1
$lm348.NumberFormat = "@"

# --- Row 349: new positive cases 75 -> 100 ---
$ws.Range("C349").Value = 100

# --- Row 350: new positive cases 12 -> 47; new hospital death 0 -> 1 ---
$ws.Range("C350").Value = 47
$lm350 = $ws.Range("L350")
$lm350.NumberFormat = "General"
$lm350.Value = 1
$lm350.NumberFormat = "@"

# --- Row 351: newly reported day, previously blank ---
$ws.Range("C351").Value = 15
$ws.Range("E351").Value = 14
$ws.Range("F351").Value = 10
$ws.Range("G351").Value = 81
$lm351 = $ws.Range("L351:M351")
$lm351.NumberFormat = "General"
$ws.Range("L351").Value = 0
$ws.Range("M351").Value = 0
$lm351.NumberFormat = "@"
